$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "country" column ---
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "country"

$ws.Range("F2:F19").Value = "India"

# --- Data corrections ---
# shalini: city noida -> bangalore
$ws.Range("D4").Value = "bangalore"

# name1: age 34 -> 24, city chennai -> bhubaneswar
$ws.Range("C11").Value = 24
$ws.Range("D11").Value = "bhubaneswar"

# name5 (first occurrence, row 15): age 55 -> 25
$ws.Range("C15").Value = 25

Write-Host "done"
